# Overview_Modeling.xlsx — "added scaling test run for weakly logReg and RF
# as well as played with cv for RF averaged"
#
# Three new per-center AUC readings (with a new "Average Val AUC" summary
# string each) are filled in on the already-present scaled-center rows
# (27 = LogReg, 28 = RF), and a brand-new row (39) is appended for an
# RF run averaged/restricted with cross-validation. The table and the
# sheet's used range both grow from A1:O38 to A1:O39 to cover it.
#
# NOTE on write order: the workbook's shared-string table appends new
# unique strings in the order cells are first written. The target file's
# new strings are 128="0.653 (0.054)" (M39), 129="0.645 (0.043)" (M28),
# 130="0.600 (0.024)" (M27) — so row 39 is populated first, then row 28,
# then row 27, to reproduce that exact ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 39: RF run, averaged, restricted + cv -------------------------
$ws.Range("B39").Value = "RandomForest"
$ws.Range("C39").Value = "MoCo"
$ws.Range("D39").Value = "Centers"
$ws.Range("E39").Value = "1 x 3"
$ws.Range("F39").Value = "average"
$ws.Range("J39").Value = 0.6
$ws.Range("K39").Value = 0.726
$ws.Range("L39").Value = 0.633
$ws.Range("M39").Value = "0.653 (0.054)"

# --- Row 28: RF scaled-center run, fill in the per-center AUC scores -------
$ws.Range("J28").Value = 0.629
$ws.Range("K28").Value = 0.703
$ws.Range("L28").Value = 0.602
$ws.Range("M28").Value = "0.645 (0.043)"

# --- Row 27: LogReg scaled-center run, fill in the per-center AUC scores ---
$ws.Range("J27").Value = 0.567
$ws.Range("K27").Value = 0.619
$ws.Range("L27").Value = 0.615
$ws.Range("M27").Value = "0.600 (0.024)"

# --- Grow the Tabelle1 table/autofilter to include the new row -------------
$ws.ListObjects.Item(1).Resize($ws.Range("A1:O39"))

# --- Match the author's final viewport/selection state ---------------------
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("M28").Select()
